$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text for the target languages (E2/F2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527

$wsZh.Range("J2").Value = "7f7f708c-33fa-4b05-9543-528a16ef3fee.md"
$wsZh.Hyperlinks.Add($wsZh.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/19cfc6c1cbea2280c0fcdeed546e81dccabb35ba/e2e/7f7f708c-33fa-4b05-9543-528a16ef3fee.md", "", "", "7f7f708c-33fa-4b05-9543-528a16ef3fee.md")
$wsZh.Columns.Item(10).ColumnWidth = 39.9618007114955

$wsZh.Range("K2").Value = "7f7f708c-33fa-4b05-9543-528a16ef3fee.eb5beac45b546f6ed723d5c177c6de1450a273f6.zh-cn.xlf"
$wsZh.Columns.Item(11).ColumnWidth = 40

$wsZh.Range("L2").Value = "2017-03-03 02:02:07"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527

$wsDe.Range("J2").Value = "7f7f708c-33fa-4b05-9543-528a16ef3fee.md"
$wsDe.Hyperlinks.Add($wsDe.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/19cfc6c1cbea2280c0fcdeed546e81dccabb35ba/e2e/7f7f708c-33fa-4b05-9543-528a16ef3fee.md", "", "", "7f7f708c-33fa-4b05-9543-528a16ef3fee.md")
$wsDe.Columns.Item(10).ColumnWidth = 39.9618007114955

$wsDe.Range("K2").Value = "7f7f708c-33fa-4b05-9543-528a16ef3fee.eb5beac45b546f6ed723d5c177c6de1450a273f6.de-de.xlf"
$wsDe.Columns.Item(11).ColumnWidth = 40

$wsDe.Range("L2").Value = "2017-03-03 02:02:28"
